# This edit reshuffles the weekly price-report rows for
# "Agrícola del Norte S.A. de Arica - Pera": each data row (2-24) keeps its
# constant identifying columns (A,B,C,E,F,G,H,I,J) but the variable columns
# (D: Fecha, K: Variedad, L: Calidad, M: Volumen, N: Precio mínimo,
# O: Precio máximo, P: Precio promedio ponderado, Q: Unidad de
# comercialización, R: Origen, S: Precio $/Kg, T: Kg/unidad) get permuted to
# a different row, matching a newer weekly snapshot of the same dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read all the source row data (columns D:T) for rows 2-24 first, before
# writing anything, so that overlapping/cyclic row moves do not clobber data
# that still needs to be read.
$rowData = @{}
$rowData["2"]  = $ws.Range("D2:T2").Value()
$rowData["3"]  = $ws.Range("D3:T3").Value()
$rowData["4"]  = $ws.Range("D4:T4").Value()
$rowData["5"]  = $ws.Range("D5:T5").Value()
$rowData["6"]  = $ws.Range("D6:T6").Value()
$rowData["7"]  = $ws.Range("D7:T7").Value()
$rowData["8"]  = $ws.Range("D8:T8").Value()
$rowData["9"]  = $ws.Range("D9:T9").Value()
$rowData["10"] = $ws.Range("D10:T10").Value()
$rowData["11"] = $ws.Range("D11:T11").Value()
$rowData["12"] = $ws.Range("D12:T12").Value()
$rowData["13"] = $ws.Range("D13:T13").Value()
$rowData["14"] = $ws.Range("D14:T14").Value()
$rowData["15"] = $ws.Range("D15:T15").Value()
$rowData["16"] = $ws.Range("D16:T16").Value()
$rowData["17"] = $ws.Range("D17:T17").Value()
$rowData["18"] = $ws.Range("D18:T18").Value()
$rowData["19"] = $ws.Range("D19:T19").Value()
$rowData["20"] = $ws.Range("D20:T20").Value()
$rowData["21"] = $ws.Range("D21:T21").Value()
$rowData["22"] = $ws.Range("D22:T22").Value()
$rowData["23"] = $ws.Range("D23:T23").Value()
$rowData["24"] = $ws.Range("D24:T24").Value()

# Map: target row -> source row whose D:T content it should now hold.
$rowMap = [ordered]@{
    "2"  = "23"
    "3"  = "24"
    "4"  = "2"
    "5"  = "3"
    "6"  = "17"
    "7"  = "18"
    "8"  = "15"
    "9"  = "16"
    "10" = "10"
    "11" = "11"
    "12" = "4"
    "13" = "5"
    "14" = "14"
    "15" = "19"
    "16" = "20"
    "17" = "21"
    "18" = "22"
    "19" = "12"
    "20" = "13"
    "21" = "6"
    "22" = "8"
    "23" = "9"
    "24" = "7"
}

foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $ws.Range("D" + $target + ":T" + $target).Value = $rowData[$source]
}
